$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the report title (was mis-typed without diacritics)
$ws.Range("A3").Value = "DANH SÁCH GIÁO VIÊN"

# Row 9 (Nguyễn Mạnh Hiếu's record) gets a real address + phone number
$ws.Range("B9").Value = "NGUYỄN VĂN HÙNG"
$ws.Range("C9").Value = "HÀ TÂY"

# Phone numbers / ids can start with a leading zero, so format as text
# before typing them in to keep the leading zero intact.
$ws.Range("D9:D11").NumberFormat = "@"

$ws.Range("D9").Value = "031231231"

# Row 10 - new teacher record (was placeholder "1","1","1")
$ws.Range("B10").Value = "HÙNG NGUYỄN VĂN"
$ws.Range("C10").Value = "ỨNG HÒA"
$ws.Range("D10").Value = "123123123"

# Row 11 - new teacher record (was placeholder "2","2","2")
$ws.Range("B11").Value = "Mạnh hiếu"
$ws.Range("C11").Value = "sóc sơn"
$ws.Range("D11").Value = "123123"
